$d = $word.ActiveDocument

# The paragraph currently reads "...has the best performnace," where
# "performnace" is wrapped in spellcheck proofErr markers and is
# followed by a separate run containing just a trailing comma. We need
# it to read "...has the best performance" (correct spelling, comma
# removed, no leftover proofErr markers), while leaving the preceding
# run ("...has the best ") completely untouched.

$rng = $d.Content
$found = $rng.Find.Execute("performnace,", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$start = $rng.Start
$rng.Text = "performance"
$wordLen = 11  # Len("performance")
$end = $start + $wordLen

# Merge the fixed run back with the run before it (they share the same
# formatting) by making, then immediately undoing, a one-character edit
# that straddles the boundary. This clears the stale proofErr markers
# that used to flank the misspelled word.
$boundary = $d.Range($start - 1, $start + 1)
$boundaryText = $boundary.Text
$boundary.Text = $boundaryText.Substring(0, 1) + "Z"
$d.Range($start - 1, $start + 1).Text = $boundaryText

# Toggling character formatting on just the corrected word forces Word
# to split it back out into its own run, giving back the original two
# run structure but now free of any proofErr markers.
$wordRng = $d.Range($start, $end)
$wordRng.Bold = 1
$wordRng.Bold = 0
